$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("tandartsenpraktijk", "Staas & Bergmans, locatie Groote Wielen", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-staas-bergmans-locatie-groote-wielen-rosmalen-10025334")
    ,@("tandartsenpraktijk", "Tandarts Abas", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandarts-abas-bussum-3054204")
    ,@("tandartsenpraktijk", "De Kliniek voor Tandheelkunde", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-de-kliniek-voor-tandheelkunde-utrecht-182964")
    ,@("tandartsenpraktijk", "Dijkstra & De Wet Tandheelkunde", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dijkstra-de-wet-tandheelkunde-amsterdam-10020311")
    ,@("tandartsenpraktijk", "Amstel Dental - Tandarts Bussum", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-amstel-dental-tandarts-bussum-bussum-10021599")
    ,@("tandartsenpraktijk", "Dental Clinics Oldenzaal", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-oldenzaal-oldenzaal-3048472")
    ,@("tandartsenpraktijk", "Venema tandartsen", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-venema-tandartsen-helmond-185417")
    ,@("tandartsenpraktijk", "KT3 tandartspraktijk", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-kt3-tandartspraktijk-zaandam-180596")
    ,@("tandartsenpraktijk", "Tandartspraktijk Romero", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-romero-honselersdijk-10027597")
    ,@("tandartsenpraktijk", "De Tandarts, kliniek voor tandheelkunde en implantologie", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-de-tandarts-kliniek-voor-tandheelkunde-en-implantologie-raamsdonksveer-184568")
    ,@("tandartsenpraktijk", "Operatandarts de By & de By", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-operatandarts-de-by-de-by-amsterdam-209685")
    ,@("tandartsenpraktijk", "Mondzorg Oosterhout", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-mondzorg-oosterhout-oosterhout-10026683")
    ,@("tandartsenpraktijk", "Tandartsenpraktijk Paul en de Witte", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartsenpraktijk-paul-en-de-witte-vlissingen-184137")
    ,@("tandartsenpraktijk", "Tandheelkunde Beverdam", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkunde-beverdam-wierden-3036333")
    ,@("tandartsenpraktijk", "Dental Clinics Zaltbommel", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-zaltbommel-zaltbommel-208738")
    ,@("tandartsenpraktijk", "Tandartsenpraktijk Grou", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartsenpraktijk-grou-grou-188160")
    ,@("tandartsenpraktijk", "Dental Clinics Gouda Greenline", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-gouda-greenline-gouda-182093")
    ,@("tandartsenpraktijk", "PC2-tandartsen", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-pc2-tandartsen-eindhoven-3032943")
    ,@("tandartsenpraktijk", "Dental Clinics Beuningen", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-beuningen-beuningen-3046007")
    ,@("tandartsenpraktijk", "Dental Clinics Colmschate", "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-colmschate-colmschate-187204")
)

$startRow = 102
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
